$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 50

# Column A holds a value that looks like a date ("2025-09-30"). Excel's
# smart-parsing would otherwise convert it into a date serial number when
# assigned via .Value, which would not match the source data (plain text,
# same as every other row in this sheet). Force the cell to Text format
# first so the string is stored verbatim, then restore the cell style to
# the sheet's normal/default style so no stray formatting is introduced.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-30"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "21:20:53"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,725.9574"
